# A new weekly Ciboulette price record was added for Femacal de La Calera.
# In the sheet, row 157 held the most-recent record (serial date 44915,
# i.e. 2022-12-20). A new row is inserted right below it, duplicating that
# record, and the newer record's date (serial 44930, i.e. 2023-01-04) is
# written into the original row 157 - pushing all the older history rows
# down by one position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 158 - this shifts rows 158:393 down to 159:394
# and extends the sheet's used range to row 394.
$ws.Rows.Item(158).Insert()

# Duplicate row 157 (the most recent record at the time) into the freshly
# inserted row 158, so the old data point is preserved one row down.
$ws.Range("A157:R157").Copy($ws.Range("A158:R158"))

# Row 157 now represents the new, more recent price entry - update its date.
$ws.Range("D157").Value = 44930
